# Applies the cryptos list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Using an ordered set of
# assignments keeps every value (including number-looking strings like
# "1.495" or percentages like "  +0.35%  ") stored as literal text,
# matching the inline-string cells already used in the sheet.
$updates = [ordered]@{
    'D2' = '24.713.07'
    'E2' = '  +0.35%  '
    'D3' = '1.694.37'
    'E3' = '  -0.14%  '
    'E4' = '  +0.32%  '
    'D5' = '317.19'
    'E5' = '  +1.34%  '
    'E6' = '  +0.44%  '
    'D7' = '0.3957'
    'E7' = '  +0.23%  '
    'E8' = '  +0.98%  '
    'D9' = '1.495'
    'E9' = '  -1.69%  '
    'E10' = '  +0.41%  '
    'D11' = '52.78'
    'E11' = '  -2.85%  '
    'D12' = '0.08933'
    'E12' = '  +2.05%  '
    'D13' = '7.274'
    'E13' = '  -0.68%  '
    'D14' = '23.64'
    'E14' = '  +2.25%  '
    'D15' = '8.078'
    'E15' = '  +6.35%  '
    'D16' = '0.00001322'
    'E16' = '  +0.18%  '
    'D17' = '1.697.75'
    'E17' = '  +0.00%  '
    'D18' = '100.05'
    'E18' = '  -0.24%  '
    'D19' = '0.07043'
    'E19' = '  -0.45%  '
    'D20' = '19.65'
    'E20' = '  +1.04%  '
    'E21' = '  +4.47%  '
    'E22' = '  +0.24%  '
    'D23' = '14.34'
    'E23' = '  +1.42%  '
    'D24' = '24.686.32'
    'E24' = '  +0.31%  '
    'D25' = '3.269'
    'E25' = '  +7.34%  '
    'D26' = '2.359'
    'E26' = '  +2.25%  '
    'E27' = '  +1.95%  '
    'D28' = '162.23'
    'E28' = '  +1.83%  '
    'B29' = 'BitcoinCash'
    'C29' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D29' = '136.13'
    'E29' = '  +2.04%  '
    'B30' = 'HuobiToken'
    'C30' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D30' = '5.200'
    'E30' = '  +0.56%  '
    'D31' = '7.515'
    'E31' = '  -1.28%  '
    'D32' = '0.08638'
    'E32' = '  +0.16%  '
    'D33' = '1.055'
    'E33' = '  -3.57%  '
    'D34' = '7.056'
    'E34' = '  -3.89%  '
    'D35' = '11.38'
    'E35' = '  +2.93%  '
    'D36' = '0.2743'
    'E36' = '  +0.85%  '
    'D37' = '1.885'
    'E37' = '  -3.93%  '
    'E38' = '  -1.78%  '
    'D39' = '0.09249'
    'E39' = '  +3.07%  '
    'D40' = '0.02729'
    'E40' = '  -0.82%  '
    'E41' = '  +0.31%  '
    'D42' = '0.7680'
    'E42' = '  +0.42%  '
    'D43' = '16.25'
    'E43' = '  +5.34%  '
    'D44' = '2.605'
    'E44' = '  +6.64%  '
    'D45' = '0.7167'
    'E45' = '  -0.14%  '
    'D46' = '4.228'
    'E46' = '  +1.38%  '
    'E47' = '  +0.44%  '
    'D48' = '140.35'
    'E48' = '  -0.11%  '
    'D49' = '1.321'
    'E49' = '  +0.07%  '
    'D50' = '90.78'
    'E50' = '  +5.04%  '
    'D51' = '0.07988'
    'E51' = '  -0.30%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text number format first so values like "1.495" or "5.200"
    # are not reinterpreted as numbers/dates by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
